$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 127
$ws.Range("F2").Value = 91
$ws.Range("H2").Value = 97
$ws.Range("E3").Value = 52
$ws.Range("F3").Value = 38
$ws.Range("H3").Value = 39
$ws.Range("E4").Value = 54
$ws.Range("F4").Value = 34
$ws.Range("H4").Value = 46
$ws.Range("E5").Value = 180
$ws.Range("F5").Value = 127
$ws.Range("H5").Value = 138
$ws.Range("E6").Value = 61
$ws.Range("F6").Value = 41
$ws.Range("H6").Value = 51
$ws.Range("F7").Value = 35
$ws.Range("H7").Value = 39
$ws.Range("E8").Value = 12
$ws.Range("E10").Value = 823
$ws.Range("F10").Value = 482
$ws.Range("H10").Value = 577
$ws.Range("E11").Value = 547
$ws.Range("F11").Value = 327
$ws.Range("H11").Value = 392
$ws.Range("E12").Value = 847
$ws.Range("F12").Value = 532
$ws.Range("H12").Value = 618
$ws.Range("E13").Value = 185
$ws.Range("F13").Value = 109
$ws.Range("H13").Value = 143
$ws.Range("F14").Value = 96
$ws.Range("H14").Value = 130
$ws.Range("E16").Value = 253
$ws.Range("E17").Value = 141
$ws.Range("E19").Value = 17
$ws.Range("F19").Value = 10
$ws.Range("H19").Value = 13
$ws.Range("E21").Value = 164
$ws.Range("F21").Value = 103
$ws.Range("H21").Value = 134
$ws.Range("E22").Value = 223
$ws.Range("E23").Value = 261
$ws.Range("F23").Value = 134
$ws.Range("H23").Value = 186
$ws.Range("E24").Value = 325
$ws.Range("F24").Value = 186
$ws.Range("H24").Value = 216
$ws.Range("E25").Value = 386
$ws.Range("F25").Value = 222
$ws.Range("H25").Value = 282
$ws.Range("E26").Value = 252
$ws.Range("F26").Value = 153
$ws.Range("H26").Value = 178
$ws.Range("E27").Value = 441
$ws.Range("F27").Value = 258
$ws.Range("H27").Value = 340
$ws.Range("E28").Value = 260
$ws.Range("F28").Value = 130
$ws.Range("H28").Value = 182
$ws.Range("E29").Value = 206
$ws.Range("F29").Value = 129
$ws.Range("H29").Value = 170
$ws.Range("E32").Value = 240
$ws.Range("F32").Value = 158
$ws.Range("H32").Value = 196
$ws.Range("E33").Value = 377
$ws.Range("E34").Value = 287
$ws.Range("F34").Value = 204
$ws.Range("H34").Value = 242
$ws.Range("E36").Value = 97
$ws.Range("F36").Value = 62
$ws.Range("H36").Value = 72
$ws.Range("E37").Value = 218
$ws.Range("F37").Value = 123
$ws.Range("H37").Value = 159
$ws.Range("E38").Value = 116
$ws.Range("E40").Value = 341
$ws.Range("F40").Value = 188
$ws.Range("H40").Value = 268
$ws.Range("E41").Value = 498
$ws.Range("F41").Value = 267
$ws.Range("H41").Value = 359
$ws.Range("E42").Value = 530
$ws.Range("F42").Value = 323
$ws.Range("H42").Value = 384
$ws.Range("E43").Value = 164
$ws.Range("F43").Value = 97
$ws.Range("H43").Value = 124
$ws.Range("E44").Value = 438
$ws.Range("F44").Value = 250
$ws.Range("H44").Value = 318
$ws.Range("E45").Value = 207
$ws.Range("F45").Value = 126
$ws.Range("H45").Value = 165
$ws.Range("E46").Value = 438
$ws.Range("F46").Value = 261
$ws.Range("H46").Value = 325
$ws.Range("E47").Value = 633
$ws.Range("F47").Value = 369
$ws.Range("H47").Value = 461
$ws.Range("E48").Value = 307
$ws.Range("F48").Value = 157
$ws.Range("H48").Value = 201
$ws.Range("F49").Value = 194
$ws.Range("H49").Value = 281
$ws.Range("F50").Value = 179
$ws.Range("H50").Value = 252
$ws.Range("E51").Value = 279
$ws.Range("F51").Value = 146
$ws.Range("H51").Value = 220
$ws.Range("E52").Value = 37